$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the test data: password value becomes a text string instead of a number
$ws.Range("B2").Value = "123456X"

# Move the active selection to D11 as recorded in the author's session
$ws.Range("D11").Select()
